$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("COVID_deaths_by_ethnicity")
$raw = $wb.Worksheets.Item("raw_data")

# Update the "data as of" caption text (C3 on the main sheet)
$main.Range("C3").Value = "Data for England up until 12th May 2020 (published 14th May 2020)"

# Refresh the underlying ONS population / NHS England death counts (raw_data sheet)
# with the updated figures as of 12th May 2020 (published 14th May 2020).
$raw.Range("C2").Value = 819402
$raw.Range("D2").Value = 0.0154567824588243
$raw.Range("E2").Value = 366
$raw.Range("F2").Value = 0.016847726017307999
$raw.Range("G2").Value = 336
$raw.Range("H2").Value = 105
$raw.Range("I2").Value = 0.0048257658015763401
$raw.Range("J2").Value = 261
$raw.Range("K2").Value = 3.48571428571428
$raw.Range("C3").Value = 436514
$raw.Range("D3").Value = 0.0082341780203505303
$raw.Range("E3").Value = 146
$raw.Range("F3").Value = 0.0067206775916037499
$raw.Range("G3").Value = 179
$raw.Range("H3").Value = 37
$raw.Range("I3").Value = 0.00169538747266972
$raw.Range("J3").Value = 109
$raw.Range("K3").Value = 3.9459459459459398
$raw.Range("C4").Value = 379503
$raw.Range("D4").Value = 0.0071587515205860297
$raw.Range("E4").Value = 78
$raw.Range("F4").Value = 0.0035904989872951498
$raw.Range("G4").Value = 156
$raw.Range("H4").Value = 53
$raw.Range("I4").Value = 0.00245594583315349
$raw.Range("J4").Value = 25
$raw.Range("K4").Value = 1.47169811320754
$raw.Range("C5").Value = 1395702
$raw.Range("D5").Value = 0.026327812467318999
$raw.Range("E5").Value = 665
$raw.Range("F5").Value = 0.0306113054686061
$raw.Range("G5").Value = 572
$raw.Range("H5").Value = 277
$raw.Range("I5").Value = 0.012745207923964399
$raw.Range("J5").Value = 388
$raw.Range("K5").Value = 2.40072202166064
$raw.Range("C6").Value = 1112282
$raw.Range("D6").Value = 0.020981521776693299
$raw.Range("E6").Value = 448
$raw.Range("F6").Value = 0.020622353157797799
$raw.Range("G6").Value = 456
$raw.Range("H6").Value = 121
$raw.Range("I6").Value = 0.0055731356410174803
$raw.Range("J6").Value = 327
$raw.Range("K6").Value = 3.70247933884297
$raw.Range("C7").Value = 977741
$raw.Range("D7").Value = 0.018443608800165701
$raw.Range("E7").Value = 394
$raw.Range("F7").Value = 0.018136623089670401
$raw.Range("G7").Value = 401
$raw.Range("H7").Value = 75
$raw.Range("I7").Value = 0.0034401856851248901
$raw.Range("J7").Value = 319
$raw.Range("K7").Value = 5.2533333333333303
$raw.Range("C8").Value = 277857
$raw.Range("D8").Value = 0.0052413530887910498
$raw.Range("E8").Value = 199
$raw.Range("F8").Value = 0.0091603756214325097
$raw.Range("G8").Value = 114
$raw.Range("H8").Value = 25
$raw.Range("I8").Value = 0.0011331290452152201
$raw.Range("J8").Value = 174
$raw.Range("K8").Value = 7.96
$raw.Range("C9").Value = 591016
$raw.Range("D9").Value = 0.011148625145758101
$raw.Range("E9").Value = 601
$raw.Range("F9").Value = 0.027665255017492101
$raw.Range("G9").Value = 242
$raw.Range("H9").Value = 183
$raw.Range("I9").Value = 0.0084047637644137797
$raw.Range("J9").Value = 418
$raw.Range("K9").Value = 3.2841530054644799
$raw.Range("C10").Value = 283005
$raw.Range("D10").Value = 0.0053384623417560501
$raw.Range("E10").Value = 69
$raw.Range("F10").Value = 0.00317621064260725
$raw.Range("G10").Value = 116
$raw.Range("H10").Value = 26
$raw.Range("I10").Value = 0.0011771409134594701
$raw.Range("J10").Value = 43
$raw.Range("K10").Value = 2.6538461538461502
$raw.Range("C11").Value = 332708
$raw.Range("D11").Value = 0.0062760344474513596
$raw.Range("E11").Value = 28
$raw.Range("F11").Value = 0.00128889707236236
$raw.Range("G11").Value = 136
$raw.Range("H11").Value = 26
$raw.Range("I11").Value = 0.00121156151643394
$raw.Range("J11").Value = 2
$raw.Range("K11").Value = 1.07692307692307
$raw.Range("C12").Value = 161550
$raw.Range("D12").Value = 0.0030473970117513502
$raw.Range("E12").Value = 14
$raw.Range("F12").Value = 0.00064444853618118197
$raw.Range("G12").Value = 66
$raw.Range("H12").Value = 9
$raw.Range("I12").Value = 0.00042803486107711102
$raw.Range("J12").Value = 5
$raw.Range("K12").Value = 1.55555555555555
$raw.Range("C13").Value = 415616
$raw.Range("D13").Value = 0.0078399687801674299
$raw.Range("E13").Value = 44
$raw.Range("F13").Value = 0.00202540968514085
$raw.Range("G13").Value = 170
$raw.Range("H13").Value = 35
$raw.Range("I13").Value = 0.0016223946394365899
$raw.Range("J13").Value = 9
$raw.Range("K13").Value = 1.25714285714285
$raw.Range("C14").Value = 548418
$raw.Range("D14").Value = 0.0103450781454079
$raw.Range("E14").Value = 585
$raw.Range("F14").Value = 0.026928742404713601
$raw.Range("G14").Value = 225
$raw.Range("H14").Value = 71
$raw.Range("I14").Value = 0.00327860075308354
$raw.Range("J14").Value = 514
$raw.Range("K14").Value = 8.23943661971831
$raw.Range("C15").Value = 2484905
$raw.Range("D15").Value = 0.046873983729408798
$raw.Range("E15").Value = 751
$raw.Range("F15").Value = 0.034570060762290501
$raw.Range("G15").Value = 1018
$raw.Range("H15").Value = 445
$raw.Range("I15").Value = 0.0204696536462134
$raw.Range("J15").Value = 306
$raw.Range("K15").Value = 1.6876404494382
$raw.Range("C16").Value = 42279236
$raw.Range("D16").Value = 0.79753399842482298
$raw.Range("E16").Value = 17117
$raw.Range("F16").Value = 0.78793039955809197
$raw.Range("G16").Value = 17326
$raw.Range("H16").Value = 19886
$raw.Range("I16").Value = 0.91539385264806705
$raw.Range("J16").Value = -2769
$raw.Range("K16").Value = 0.86075631097254302
$raw.Range("C17").Value = 517001
$raw.Range("D17").Value = 0.0097524438407456497
$raw.Range("E17").Value = 219
$raw.Range("F17").Value = 0.010081016387405601
$raw.Range("G17").Value = 212
$raw.Range("H17").Value = 351
$raw.Range("I17").Value = 0.0161452398550926
$raw.Range("J17").Value = -132
$raw.Range("K17").Value = 0.62393162393162305

# Recalculate all dependent formulas on the summary sheet
$wb.Application.Calculate()

# Restore the active selection on the main sheet
$main.Activate()
$main.Range("C4").Select()
